{"js": "// 1) Intro paragraph: the four runs that together spell out the sentence\n//    \"Banknotes are one of the most important assets ... genuine or not.\"\n//    get merged into a single run. Re-inserting the identical text over\n//    the found range collapses the run boundaries without altering the\n//    wording.\nconst introText =\n  \"Banknotes are one of the most important assets of a country. Some \" +\n  \"miscreants introduce fake notes which bear a resemblance to original \" +\n  \"note to create discrepancies of the money in the financial market. It \" +\n  \"is difficult for humans to tell true and fake banknotes apart \" +\n  \"especially because they have a lot of similar features. Fake notes \" +\n  \"are created with precision, hence there is need for an efficient \" +\n  \"algorithm (ANN) which accurately predicts whether a banknote is \" +\n  \"genuine or not.\";\n\nconst introSearch = context.document.body.search(introText, { matchCase: true });\nintroSearch.load(\"text\");\nawait context.sync();\n\nif (introSearch.items.length > 0) {\n  introSearch.items[0].insertText(introText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) CONCLUSION paragraph: trim the opening clause down to \"ANN techniques\n//    used to detect forged banknotes.\" and drop \"support vector machine\"\n//    from the result sentence.\nconst openingOld =\n  \"various techniques used to detect forged banknotes, this paper \" +\n  \"presents banknote authentication for recognizing the banknote as \" +\n  \"genuine or fake by using two supervised learning techniques. \" +\n  \"Extensive experiments have been performed on banknotes dataset using \" +\n  \"both the models to find the best model suitable for classification \" +\n  \"of the notes.\";\nconst openingNew = \"ANN techniques used to detect forged banknotes.\";\n\nconst openingSearch = context.document.body.search(openingOld, { matchCase: true });\nopeningSearch.load(\"text\");\nawait context.sync();\n\nif (openingSearch.items.length > 0) {\n  openingSearch.items[0].insertText(openingNew, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst svmOld = \"outperforms support vector machine and gives\";\nconst svmNew = \"outperforms and gives\";\n\nconst svmSearch = context.document.body.search(svmOld, { matchCase: true });\nsvmSearch.load(\"text\");\nawait context.sync();\n\nif (svmSearch.items.length > 0) {\n  svmSearch.items[0].insertText(svmNew, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Intro paragraph: the sentence about banknotes/ANN is currently split\n#    across four runs. Re-running Find & Replace over the identical text\n#    collapses it back down to a single run, the same way Word merges\n#    runs whenever it rewrites a matched range.\n$introText = \"Banknotes are one of the most important assets of a country. Some miscreants introduce fake notes which bear a resemblance to original note to create discrepancies of the money in the financial market. It is difficult for humans to tell true and fake banknotes apart especially because they have a lot of similar features. Fake notes are created with precision, hence there is need for an efficient algorithm (ANN) which accurately predicts whether a banknote is genuine or not.\"\n\n$introRange = $d.Content\n$introFind = $introRange.Find\n$introFind.ClearFormatting()\n$introFind.Replacement.ClearFormatting()\n$introFind.Text = $introText\n$introFind.Replacement.Text = $introText\n$introFind.MatchCase = $true\n$introFind.MatchWildcards = $false\n$introFind.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 2) CONCLUSION paragraph: shorten the opening clause to\n#    \"ANN techniques used to detect forged banknotes.\"\n$openingRange = $d.Content\n$openingFind = $openingRange.Find\n$openingFind.ClearFormatting()\n$openingFind.Replacement.ClearFormatting()\n$openingFind.Text = \"various techniques used to detect forged banknotes, this paper presents banknote authentication for recognizing the banknote as genuine or fake by using two supervised learning techniques. Extensive experiments have been performed on banknotes dataset using both the models to find the best model suitable for classification of the notes.\"\n$openingFind.Replacement.Text = \"ANN techniques used to detect forged banknotes.\"\n$openingFind.MatchCase = $true\n$openingFind.MatchWildcards = $false\n$openingFind.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 3) Drop \"support vector machine\" from the result sentence.\n$svmRange = $d.Content\n$svmFind = $svmRange.Find\n$svmFind.ClearFormatting()\n$svmFind.Replacement.ClearFormatting()\n$svmFind.Text = \"outperforms support vector machine and gives\"\n$svmFind.Replacement.Text = \"outperforms and gives\"\n$svmFind.MatchCase = $true\n$svmFind.MatchWildcards = $false\n$svmFind.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n"}
